# Weekly update: insert a new Ajo (garlic) price observation as the new
# row 120 on the sheet, pushing the existing rows 120-136 down to 121-137.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 120 (shifts rows 120:136 -> 121:137)
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with this week's data point
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 44491
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112003
$ws.Cells.Item(120, 7).Value = "Ajo"
$ws.Cells.Item(120, 8).Value = "Chino"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 60
$ws.Cells.Item(120, 11).Value = 16000
$ws.Cells.Item(120, 12).Value = 17000
$ws.Cells.Item(120, 13).Value = 16500
$ws.Cells.Item(120, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(120, 15).Value = "China"
$ws.Cells.Item(120, 16).Value = 1650
$ws.Cells.Item(120, 17).Value = 10
$ws.Cells.Item(120, 18).Value = "Hortaliza"
